$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.858.63'
$ws.Range("E2").Value = '  -0.04%  '

$ws.Range("D3").Value = '1.893.25'
$ws.Range("E3").Value = '  -0.26%  '

$ws.Range("E4").Value = '  -0.44%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7766'
$ws.Range("E5").Value = '  -2.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.85'
$ws.Range("E6").Value = '  +0.87%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.47%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3146'
$ws.Range("E8").Value = '  -1.77%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07403'
$ws.Range("E9").Value = '  +4.22%  '

$ws.Range("E10").Value = '  -3.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08137'
$ws.Range("E11").Value = '  +1.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7673'
$ws.Range("E12").Value = '  -0.40%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.497'
$ws.Range("E13").Value = '  +3.50%  '

$ws.Range("D14").Value = '1.863.84'
$ws.Range("E14").Value = '  -5.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.49'
$ws.Range("E15").Value = '  -0.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.235'
$ws.Range("E16").Value = '  +5.37%  '

$ws.Range("D17").Value = '29.872.58'
$ws.Range("E17").Value = '  -0.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.99'
$ws.Range("E18").Value = '  +0.77%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.31'
$ws.Range("E19").Value = '  +0.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007860'
$ws.Range("E20").Value = '  +1.59%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.46%  '

$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.168'
$ws.Range("E22").Value = '  +0.30%  '

$ws.Range("D23").Value = '2.143.40'
$ws.Range("E23").Value = '  -2.56%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  -0.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1573'
$ws.Range("E25").Value = '  -2.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.439'
$ws.Range("E26").Value = '  +1.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.39'
$ws.Range("E27").Value = '  -2.00%  '

$ws.Range("E28").Value = '  +0.59%  '

$ws.Range("E29").Value = '  -2.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.449'
$ws.Range("E30").Value = '  +5.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.550'
$ws.Range("E31").Value = '  +0.79%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.518'
$ws.Range("E32").Value = '  +0.89%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05633'
$ws.Range("E33").Value = '  -0.84%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.105'
$ws.Range("E34").Value = '  +0.61%  '

$ws.Range("E35").Value = '  -1.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7686'
$ws.Range("E36").Value = '  +4.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.003'
$ws.Range("E37").Value = '  +0.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.644'
$ws.Range("E38").Value = '  -2.72%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01931'
$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.790'
$ws.Range("E40").Value = '  +0.40%  '

$ws.Range("D41").Value = '1.160.57'
$ws.Range("E41").Value = '  +13.53%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4476'
$ws.Range("E42").Value = '  +0.45%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '74.35'
$ws.Range("E43").Value = '  +2.92%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.997'
$ws.Range("E44").Value = '  +1.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8551'
$ws.Range("E45").Value = '  +0.99%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.911'
$ws.Range("E46").Value = '  +1.23%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  -0.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.19'
$ws.Range("E48").Value = '  -0.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.935'
$ws.Range("E49").Value = '  +1.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.099'
$ws.Range("E50").Value = '  +2.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.560'
$ws.Range("E51").Value = '  +0.91%  '
